$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for number-looking values (prices & some rating counts)
# so Excel keeps the literal text (e.g. "199.", "70  ") instead of coercing to a number
$ws.Range("B2:B16").NumberFormat = "@"
$ws.Range("C5:C7").NumberFormat = "@"
$ws.Range("C11:C13").NumberFormat = "@"
$ws.Range("C17").NumberFormat = "@"

# Populate cell values (scraped Amazon phone listings: title/price/rating)
$ws.Range('A1').Value = 'title'
$ws.Range('B1').Value = 'price'
$ws.Range('C1').Value = 'rating'
$ws.Range('A2').Value = 'OUKITEL WP36 2024 Rugged Smartphone - 10600mAh Rugged Phone 128dB Loud Speaker, 16GB+128GB Dual Sim Phones 6.52" Big Screen, Android 13 Cell Phone, NFC/OTG, T-mobile Compatible Rugged Smart Phone  '
$ws.Range('B2').Value = '199.'
$ws.Range('C2').Value = '4.4 out of 5 stars '
$ws.Range('A3').Value = 'Alcatel 1 (32GB) 5.0" Full View Display - Removable Battery - Dual SIM GSM Unlocked US & Global 4G LTE International Version - Al Aqua  '
$ws.Range('B3').Value = '35.'
$ws.Range('C3').Value = '4.4 out of 5 stars'
$ws.Range('A4').Value = 'Total by Verizon TCL 30 Z, 32GB, Black - Prepaid Smartphone (Locked)  '
$ws.Range('B4').Value = '39.'
$ws.Range('C4').Value = '4.4 out of 5 stars'
$ws.Range('A5').Value = 'BLU G73 | 2023 | 3-Day Battery | Unlocked | 6.8” HD+ Infinity Display | 128/6GB | Triple 50MP Camera | US Version | US Warranty | Grey  '
$ws.Range('B5').Value = '99.'
$ws.Range('C5').Value = '70  '
$ws.Range('A6').Value = 'Samsung Galaxy A15 5G (SM-156M/DSN), 256GB 8GB RAM, Dual SIM, Factory Unlocked GSM, International Version (Wall Charger Bundle) (Light Blue)  '
$ws.Range('B6').Value = '179.'
$ws.Range('C6').Value = '70  '
$ws.Range('A7').Value = 'SAMSUNG Galaxy A25 5G A Series Cell Phone, 128GB Unlocked Android Smartphone, AMOLED Display, Advanced Triple Camera System, Expandable Storage, Power Sound w/Stereo Speakers, US Version, 2024, Black  '
$ws.Range('B7').Value = '300.'
$ws.Range('C7').Value = '70'
$ws.Range('A8').Value = 'Nokia G100 | Verizon, T-Mobile, AT&T | Android 12 | Unlocked Smartphone | 3-Day Battery | US Version | 4/128GB | 6.52-Inch Screen | 13MP Triple Camera | Polar Night  '
$ws.Range('B8').Value = '99.'
$ws.Range('C8').Value = '3.6 out of 5 stars '
$ws.Range('A9').Value = 'SAMSUNG Galaxy A05s (128GB, 4GB) 6.7" Dual SIM GSM Unlocked Global 4G LTE A057M/DS (Black)  '
$ws.Range('B9').Value = '119.'
$ws.Range('C9').Value = '3.6 out of 5 stars'
$ws.Range('A10').Value = 'SAMSUNG Galaxy A35 5G A Series Cell Phone, 128GB Unlocked Android Smartphone, AMOLED Display, Advanced Triple Camera System, Expandable Storage, Rugged Design, US Version, 2024, Awesome Lilac  '
$ws.Range('B10').Value = '359.'
$ws.Range('C10').Value = '3.6 out of 5 stars'
$ws.Range('A11').Value = 'Alcatel 1 (32GB) 5.0" Full View Display - Removable Battery - Dual SIM GSM Unlocked US & Global 4G LTE International Version - Volcano Black  '
$ws.Range('B11').Value = '35.'
$ws.Range('C11').Value = '44  '
$ws.Range('A12').Value = 'Panasonic Compact Cordless Phone with DECT 6.0, 1.6" Amber LCD and Illuminated HS Keypad, Call Block, Caller ID, Multiple Display Languages - 1 Handset - KX-TGB810S (Black/Silver)  '
$ws.Range('B12').Value = '15.'
$ws.Range('C12').Value = '44  '
$ws.Range('A13').Value = 'AT&T 210 Basic Trimline Corded Phone, No AC Power Required, Wall-Mountable, White  '
$ws.Range('B13').Value = '49.'
$ws.Range('C13').Value = '44'
$ws.Range('A14').Value = 'TracFone Samsung Galaxy A03s, 32GB, Black - Prepaid Smartphone (Locked)  '
$ws.Range('B14').Value = '139.'
$ws.Range('C14').Value = '3.9 out of 5 stars '
$ws.Range('A15').Value = 'Samsung Galaxy A15 (SM-155M/DSN), 128GB 6GB RAM, Dual SIM, Factory Unlocked GSM, International Version (Ring Grip Case Bundle) (Light Blue)  '
$ws.Range('B15').Value = '569.'
$ws.Range('C15').Value = '3.9 out of 5 stars'
$ws.Range('A16').Value = '8849 Tank 3 Rugged Smartphone, 23800mAh 5G Outdoor Rugged Cell Phone Unlocked, 32GB RAM+512GB ROM, 6.79" Waterproof Android 13 Mobile Phones, 200MP Main Camera/OTG/NFC(Support T-Mobile & Verizon Only)  '
$ws.Range('B16').Value = '52.'
$ws.Range('C16').Value = '3.9 out of 5 stars'
$ws.Range('A17').Value = 'Cheap Smartphone， 5.0'''' Android 9.0, 16GB ROM(Extendable to 128GB,Dual SIM Dual Camera, WiFi,Bluetooth,GPS Basic Mobile Phones (R10-Purple)  '
$ws.Range('C17').Value = '26  '

# Restore normal (General) cell style now that the text values are locked in
$ws.Range("B2:B16").Style = "Normal"
$ws.Range("C5:C7").Style = "Normal"
$ws.Range("C11:C13").Style = "Normal"
$ws.Range("C17").Style = "Normal"

# Match row heights used by Excel when this data was entered
$ws.Range("A1:C17").RowHeight = 15.9

# Restore the last active selection
$ws.Range("H16").Select() | Out-Null
